$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 9; $row -le 45; $row++) {
    $src = $ws.Range("AD$row")
    $dst = $ws.Range("AF$row")
    $src.Copy($dst)
}
